# Revert "1.07 with intent update"
# Restores the "快递柜 / 小苏" scenario text for the 业委会工作 row (row 9)
# in place of the "徐阿姨 / 选举" scenario text, and moves the active
# selection from C10 back to C9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bg = "某商品房小区，年轻人入住率较高，快递需求量较大，但是由于工作原因无法本人签收，家里又没有老人帮忙代收，年轻人小苏希望小区引进快递柜。但是由于小区没有业委会，没有人对接快递柜的引进工作，小苏对此很不满意，于是致电12345投诉。`n接下来我来扮演小苏，而您扮演接听我电话的工作人员。"

$persona = "你叫小苏，今年30岁，你和你丈夫都是企业白领，没有孩子，平时很少与居委会有交集。最近你看到其他小区都引入了快递柜而你们小区没有，后来你得知原来是你们小区没有业委会，导致无人对接此事。你又去找居委会，居委会给出了业委会候选人名单，你强烈怀疑候选人的工作能力，认为他们也搞不定此事。你回到家后致电12345进行投诉。"

$opener = "居委会不作为啊，人家小区都有快递柜的，为什么就我们小区没有？临汾小区还是先进单位了，连个快递柜也没有，怎么服务我们年轻人啊？强烈呼吁安装快递柜！！！"

$ws.Range("B9").Value = $bg
$ws.Range("C9").Value = $persona
$ws.Range("D9").Value = $opener

$null = $ws.Range("C9").Select()
